$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 361
$ws.Cells.Item(4, 6).Value = 386
$ws.Cells.Item(5, 6).Value = 8426
$ws.Cells.Item(8, 6).Value = 2328
$ws.Cells.Item(11, 6).Value = 7716
$ws.Cells.Item(13, 6).Value = 625
$ws.Cells.Item(14, 6).Value = 155
$ws.Cells.Item(16, 6).Value = 999
$ws.Cells.Item(17, 6).Value = 1533
$ws.Cells.Item(18, 6).Value = 2159
$ws.Cells.Item(19, 6).Value = 34
$ws.Cells.Item(20, 6).Value = 220
$ws.Cells.Item(21, 6).Value = 281
$ws.Cells.Item(22, 6).Value = 1111
$ws.Cells.Item(24, 6).Value = 777
$ws.Cells.Item(25, 6).Value = 64
$ws.Cells.Item(26, 6).Value = 804
$ws.Cells.Item(27, 6).Value = 1339
$ws.Cells.Item(28, 6).Value = 523
$ws.Cells.Item(29, 6).Value = 268
$ws.Cells.Item(30, 6).Value = 24
$ws.Cells.Item(31, 6).Value = 250
$ws.Cells.Item(33, 6).Value = 71
$ws.Cells.Item(35, 6).Value = 2526

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 7856
$ws.Cells.Item(4, 6).Value = 134
$ws.Cells.Item(8, 6).Value = 24
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(18, 6).Value = 49
$ws.Cells.Item(22, 6).Value = 135
$ws.Cells.Item(29, 6).Value = 19
$ws.Cells.Item(31, 6).Value = 90
$ws.Cells.Item(41, 6).Value = 169
$ws.Cells.Item(45, 6).Value = 58

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 2434
$ws.Cells.Item(8, 6).Value = 2446
$ws.Cells.Item(9, 6).Value = 9504
$ws.Cells.Item(11, 6).Value = 199
$ws.Cells.Item(15, 6).Value = 321
$ws.Cells.Item(16, 6).Value = 2641
$ws.Cells.Item(17, 6).Value = 305
$ws.Cells.Item(18, 6).Value = 134
$ws.Cells.Item(19, 6).Value = 588
$ws.Cells.Item(8, 7).Value = 0

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 361
$ws.Cells.Item(4, 6).Value = 2434
$ws.Cells.Item(7, 6).Value = 199
$ws.Cells.Item(9, 6).Value = 321
$ws.Cells.Item(10, 6).Value = 2641
$ws.Cells.Item(11, 6).Value = 305
$ws.Cells.Item(13, 6).Value = 7716
$ws.Cells.Item(14, 6).Value = 625
$ws.Cells.Item(15, 6).Value = 155
$ws.Cells.Item(16, 6).Value = 134
$ws.Cells.Item(17, 6).Value = 999
$ws.Cells.Item(18, 6).Value = 1533
$ws.Cells.Item(19, 6).Value = 2159
$ws.Cells.Item(20, 6).Value = 134
$ws.Cells.Item(21, 6).Value = 588
$ws.Cells.Item(22, 6).Value = 588
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(26, 6).Value = 281
$ws.Cells.Item(27, 6).Value = 777
$ws.Cells.Item(28, 6).Value = 64
$ws.Cells.Item(29, 6).Value = 804
$ws.Cells.Item(31, 6).Value = 1339
$ws.Cells.Item(35, 6).Value = 523
$ws.Cells.Item(39, 6).Value = 250
$ws.Cells.Item(40, 6).Value = 19
$ws.Cells.Item(46, 6).Value = 2526
$ws.Cells.Item(49, 6).Value = 58
